# Add a new game-log entry (row 72) to the "data" sheet — script edits mm
# info in tts: add vp tag / card name, plus this extra row of results,
# and syncs the row count the shiny app expects.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("A72").Value = 2
$ws.Range("B72").Value = "The Contest of Champions"
$ws.Range("C72").Value = "Nimrod, Super Sentinel"
$ws.Range("D72").Value = "Sentinel Territories|Poisons"
$ws.Range("E72").Value = "Shi'ar Death Commandos"
$ws.Range("F72").Value = "Ant-Man (AM)|Beta Ray Bill (HOA)|Darkhawk (R)|Phyla-Vell (ITC)|Juggernaut (V)|Man-Thing (D)"
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = "26|41"
$ws.Range("I72").Value = "no"
$ws.Range("K72").Value = "Juggernaut dominated early, darkhawk had great synergy against Nimrod at the end."

# Leave the cursor where Excel would land after typing the last cell of the
# new row (matches the workbook's saved selection state).
$ws.Range("K73").Select() | Out-Null
